{"js": "// Clean up the duplicated \"egXML\" / \"egXMLTable\" paragraph styles that\n// had accumulated under auto-numbered styleIds (egXML0..egXML4,\n// egXMLTable0..egXMLTable4) because they were emitted more than once by\n// the gloss-list handling. Collapse them back down to a single egXML\n// and a single egXMLTable style definition.\n\nconst doc = context.document;\n\n// styleIds that are pure duplicates of the (identical) egXML / egXMLTable\n// paragraph style definition and can simply be dropped. Deleted in\n// reverse document order (last style first) to keep internal indices\n// stable while the collection shrinks.\nconst duplicateStyleIds = [\n  \"egXMLTable4\",\n  \"egXML4\",\n  \"egXMLTable3\",\n  \"egXML3\",\n  \"egXML2\",\n  \"egXMLTable1\",\n  \"egXMLTable0\",\n  \"egXML0\",\n  \"egXML\",\n  \"egXMLTable\"\n];\n\nfor (const styleId of duplicateStyleIds) {\n  doc.getStyles().getByName(styleId).delete();\n  await context.sync();\n}\n\n// The two surviving duplicates (egXML1 / egXMLTable2) become the new\n// canonical egXML / egXMLTable styles. The API cannot rename a style's\n// internal styleId in place, so delete them and recreate styles with\n// the canonical ids and the same formatting. (Deleted in reverse\n// document order too, same as above.)\ndoc.getStyles().getByName(\"egXMLTable2\").delete();\nawait context.sync();\ndoc.getStyles().getByName(\"egXML1\").delete();\nawait context.sync();\n\ndoc.addStyle(\"egXML\", Word.StyleType.paragraph);\nawait context.sync();\ndoc.addStyle(\"egXMLTable\", Word.StyleType.paragraph);\nawait context.sync();\n\nconst egXML = doc.getStyles().getByName(\"egXML\");\negXML.baseStyle = \"Normal\";\negXML.quickStyle = true;\negXML.font.name = \"Courier\";\negXML.font.size = 10;\nawait context.sync();\n\nconst egXMLTable = doc.getStyles().getByName(\"egXMLTable\");\negXMLTable.baseStyle = \"Normal\";\negXMLTable.quickStyle = true;\negXMLTable.font.name = \"Courier\";\negXMLTable.font.size = 9;\negXMLTable.paragraphFormat.spaceBefore = 4;\nawait context.sync();\n", "ps1": "# Clean up the duplicated \"egXML\" / \"egXMLTable\" paragraph styles that\n# had accumulated under auto-numbered styleIds (egXML0..egXML4,\n# egXMLTable0..egXMLTable4) because they were emitted more than once by\n# the gloss-list handling. Collapse them back down to a single egXML\n# and a single egXMLTable style definition.\n\n$d = $word.ActiveDocument\n\n# styleIds that are pure duplicates of the (identical) egXML / egXMLTable\n# paragraph style definition and can simply be dropped. Deleted in\n# reverse document order (last style first) to keep internal indices\n# stable while the collection shrinks.\n$duplicateStyleIds = @(\n    \"egXMLTable4\",\n    \"egXML4\",\n    \"egXMLTable3\",\n    \"egXML3\",\n    \"egXML2\",\n    \"egXMLTable1\",\n    \"egXMLTable0\",\n    \"egXML0\",\n    \"egXML\",\n    \"egXMLTable\"\n)\n\nforeach ($styleId in $duplicateStyleIds) {\n    $d.Styles($styleId).Delete()\n}\n\n# The two surviving duplicates (egXML1 / egXMLTable2) become the new\n# canonical egXML / egXMLTable styles. COM cannot rename a style's\n# internal styleId in place, so delete them and recreate styles with\n# the canonical ids and the same formatting. (Deleted in reverse\n# document order too, same as above.)\n$d.Styles(\"egXMLTable2\").Delete()\n$d.Styles(\"egXML1\").Delete()\n\n$egXML = $d.Styles.Add(\"egXML\", 1)\n$egXML.BaseStyle = \"Normal\"\n$egXML.QuickStyle = $true\n$egXML.Font.Name = \"Courier\"\n$egXML.Font.Size = 10\n\n$egXMLTable = $d.Styles.Add(\"egXMLTable\", 1)\n$egXMLTable.BaseStyle = \"Normal\"\n$egXMLTable.QuickStyle = $true\n$egXMLTable.Font.Name = \"Courier\"\n$egXMLTable.Font.Size = 9\n$egXMLTable.ParagraphFormat.SpaceBefore = 4\n"}
